# Duplicate the language/count data rows (A2:A32) below the existing
# table so they repeat starting at A33:A63 (header row A1 is not repeated).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A32").Copy()
$ws.Range("A33").PasteSpecial()
